$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character used in one PEPE price cell (U+2083)
$sub3 = [char]0x2083

$ws.Range("D2").Value = "60.005.30"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "2.494.89"
$ws.Range("E3").Value = "  -5.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.44"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.33"
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "2.519.91"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.42"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "2.938.78"
$ws.Range("E14").Value = "  -4.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.70"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "59.926.56"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "2.512.20"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.51"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.38"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.83"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.53"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.450"
$ws.Range("E25").Value = "  -9.74%  "
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.31"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "0.0$($sub3)0796"
$ws.Range("E30").Value = "  -4.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.90"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.86"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.94"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("E37").Value = "  -4.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.75"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.10"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "317.11"
$ws.Range("E40").Value = "  -6.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.74"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("E43").Value = "  -7.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.80"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.50"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0535"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0945"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0233"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("E51").Value = "  -5.29%  "
